$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assets")

# Row 3 is renamed from "TCDR_ReportMailSubject" to "TCDR_RequestMailSubject"
$ws.Range("A3").Value = "TCDR_RequestMailSubject"
$ws.Range("B3").Value = "TCDR_RequestMailSubject"
$ws.Range("C3").Value = "Logistics/Costco Trafficking"

# New row 4: the old "TCDR_ReportMailSubject" asset is re-added under its new name
$ws.Range("A4").Value = "TCDR_MailReportTransactionSubject"
$ws.Range("B4").Value = "TCDR_MailReportTransactionSubject"
$ws.Range("C4").Value = "Logistics/Costco Trafficking"

# New row 5
$ws.Range("A5").Value = "TCDR_CotscoFormOutputPath"
$ws.Range("B5").Value = "TCDR_CotscoFormOutputPath"
$ws.Range("C5").Value = "Logistics/Costco Trafficking"

# New row 6
$ws.Range("A6").Value = "TCDR_ReportPath"
$ws.Range("B6").Value = "TCDR_ReportPath"
$ws.Range("C6").Value = "Logistics/Costco Trafficking"

# Columns A:B widened to fit the longer asset names
$ws.Range("A1:B1000").ColumnWidth = 33.1667

# Update the active selection to match where the author left off editing
[void]$ws.Range("D8").Select()

Write-Host "Assets sheet updated"
